$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16").Value = "FT23168119X2K9P6"
$ws.Range("A17").Value = "FT231680RGF8G6S1"
